$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSubject")

# Update the header cell text: "Maximum contraction [Nm] / Fatigue [%]" -> "Maximum contraction [N] / Fatigue [%]"
$ws.Range("E6").Value = "Maximum contraction [N] `n/ Fatigue [%] "

# Move the active selection to F11 (matches the saved cursor position in the source file)
$ws.Range("F11").Select() | Out-Null
